$wb = $excel.ActiveWorkbook

# Rename "Shelf Sensor" sheet to "Sensor List"
$ws = $wb.Worksheets.Item("Shelf Sensor")
$ws.Name = "Sensor List"

# The old column B had wPres0, wPres1, wPres2, wTemp0, wTemp1, aFlow, aCO2, aRH, aTemp.
# The trailing general "aFlow"/"aCO2"/"aRH"/"aTemp" placeholders (rows 7-10) are being
# replaced by a dedicated, fully expanded "general_sensor" list in column C, so clear them.
$ws.Range("B7:B10").ClearContents()

# New column C: "general_sensor" header + the expanded list of general sensor tags.
$ws.Cells.Item(1, 3).Value = "general_sensor"

$generalSensors = @("wFlow0","wFlow1","aFlow0","aFlow1","aFlow2","aFlow3","aFlow4","aFlow5","aFlow6","aFlow7","aFlow8","aFlow9","aFlow10","aFlow11","aFlow12","aFlow13","aFlow14","aFlow15","aFlow16","aFlow17","aFlow18","aFlow19","aFlow20","aFlow21","aFlow22","aFlow23","aFlow24","aTemp0","aTemp1","aTemp2","aTemp3","aTemp4","aTemp5","aTemp6","aTemp7","aTemp8","aTemp9","aTemp10","aTemp11","aTemp12","aTemp13","aTemp14","aTemp15","aTemp16","aTemp17","aTemp18","aTemp19","aTemp20","aTemp21","aTemp22","aTemp23","aTemp24","aCO2_0","aCO2_1","aCO2_2","aCO2_3","aCO2_4","aCO2_5","aCO2_6","aCO2_7","aCO2_8","aCO2_9","aCO2_10","aCO2_11","aCO2_12","aCO2_13","aCO2_14","aCO2_15","aCO2_16","aCO2_17","aCO2_18","aCO2_19","aCO2_20","aCO2_21","aCO2_22","aCO2_23","aCO2_24","aRH_0","aRH_1","aRH_2","aRH_3","aRH_4","aRH_5","aRH_6","aRH_7","aRH_8","aRH_9","aRH_10","aRH_11","aRH_12","aRH_13","aRH_14","aRH_15","aRH_16","aRH_17","aRH_18","aRH_19","aRH_20","aRH_21","aRH_22","aRH_23","aRH_24")
for ($i = 0; $i -lt $generalSensors.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $generalSensors[$i]
}

# Match the page setup Excel stamps on this sheet once it has been viewed/printed.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Reflect the author's final selection on the sheet.
$ws.Range("D11").Select() | Out-Null
